$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# --- Add the new worksheet "2000-09" after the existing sheet ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "2000-09"

# --- Column widths on the new sheet ---
$ws2.Columns.Item(1).ColumnWidth = 11.25
$ws2.Columns.Item(2).ColumnWidth = 25.59

# --- Row heights ---
$ws2.Rows.Item(1).RowHeight = 129.6
$ws2.Rows.Item(5).RowHeight = 129.6

# --- Shared strings must be created in this exact order so the new ---
# --- entries land at the same shared-string indices as the target. ---
$ws2.Range("C2").Value2 = "2000-09"
$ws2.Range("B2").Value2 = "Baseline_2000-09_C81"
$ws2.Range("B3").Value2 = "Baseline_2000-09_newWeather"
$ws2.Range("A3").Value2 = "CW3M C165"

# --- Header row (row 1) ---
$ws2.Range("A1").Value2 = "model"
$ws2.Range("C1").Value2 = "Year"
$ws2.Range("D1").Value2 = " tot in HRUs reaches and reservoirs at end of last year (mm H2O)"
$ws2.Range("E1").Value2 = " Precip (mm H2O)"
$ws2.Range("F1").Value2 = " GW pumping (mm H2O)"
$ws2.Range("G1").Value2 = " High Cascades groundwater contribution mm H2O"
$ws2.Range("H1").Value2 = " from outside the basin (mm H2O)"
$ws2.Range("I1").Value2 = " water added by FlowModel (mm)"
$ws2.Range("J1").Value2 = " to outside the basin (mm H2O)"
$ws2.Range("K1").Value2 = " AET (mm H2O)"
$ws2.Range("L1").Value2 = " SNOW_EVAP (mm H2O)"
$ws2.Range("M1").Value2 = " basin discharge (mm H2O)"
$ws2.Range("N1").Value2 = " tot in HRUs reaches and reservoirs at end of this year (mm H2O)"
$ws2.Range("O1").Value2 = " irrigation (ac-ft)"
$ws2.Range("P1").Value2 = " municipal and rural domestic (ac-ft)"
$ws2.Range("Q1").Value2 = " mass balance discrepancy (mm H2O)"
$ws2.Range("R1").Value2 = " mass balance discrepancy (fraction)"
$ws2.Range("S1").Value2 = " weather year"

# Apply the header formatting (wrap text for the whole row, then per-column number formats)
$ws2.Range("A1:S1").WrapText = $true
$ws2.Range("D1:N1").NumberFormat = "0.00"
$ws2.Range("O1:P1").NumberFormat = "0"
$ws2.Range("Q1").NumberFormat = "0.00"
$ws2.Range("R1").NumberFormat = "0.000000"

# --- Row 2 data ---
$ws2.Range("A2").Value2 = "CW3M"
$ws2.Range("D2").Value2 = 1092.1001221000001
$ws2.Range("E2").Value2 = 1790.8143431000001
$ws2.Range("F2").Value2 = 1.0697084000000001
$ws2.Range("G2").Value2 = 327.85034159999998
$ws2.Range("H2").Value2 = 9.3183378000000001
$ws2.Range("J2").Value2 = 7.7646284999999988
$ws2.Range("K2").Value2 = 741.50434550000011
$ws2.Range("L2").Value2 = 85.286103100000005
$ws2.Range("M2").Value2 = 1294.5123962
$ws2.Range("N2").Value2 = 1097.7140259
$ws2.Range("O2").Value2 = 5839.5391357999997
$ws2.Range("P2").Value2 = 25979.647461100001
$ws2.Range("Q2").Value2 = 5.6286456000000005
$ws2.Range("R2").Value2 = 0.0017227000000000002
$ws2.Range("S2").Value2 = "2000-09"

# Number formats for row 2 (I2 stays blank but still needs the "0.00" style)
$ws2.Range("D2:N2").NumberFormat = "0.00"
$ws2.Range("O2:P2").NumberFormat = "0"
$ws2.Range("Q2").NumberFormat = "0.00"
$ws2.Range("R2").NumberFormat = "0.000000"

# --- Row 3 data ---
$ws2.Range("C3").Value2 = "2000-09"
$ws2.Range("D3").Value2 = 641.53857870000013
$ws2.Range("E3").Value2 = 1609.2949586000002
$ws2.Range("F3").Value2 = 1.0377343999999999
$ws2.Range("G3").Value2 = 280.39512939999997
$ws2.Range("H3").Value2 = 9.3183378000000001
$ws2.Range("I3").Value2 = 11.447250100000002
$ws2.Range("J3").Value2 = 7.7646284999999988
$ws2.Range("K3").Value2 = 634.34777839999992
$ws2.Range("L3").Value2 = 61.771183299999997
$ws2.Range("M3").Value2 = 1148.1819335
$ws2.Range("N3").Value2 = 700.97914120000007
$ws2.Range("O3").Value2 = 6439.1138917000007
$ws2.Range("P3").Value2 = 25979.647461100001
$ws2.Range("Q3").Value2 = 0.012676099999999945
$ws2.Range("R3").Value2 = -0.000038799999999999994
$ws2.Range("S3").Value2 = "2000-09"

# Number formats for row 3
$ws2.Range("D3:N3").NumberFormat = "0.00"
$ws2.Range("O3:P3").NumberFormat = "0"
$ws2.Range("Q3").NumberFormat = "0.00"
$ws2.Range("R3").NumberFormat = "0.000000"

# E3 carries the highlighted (yellow-fill) version of the "0.00" style
$ws2.Range("E3").Interior.Color = 65535

# --- Row 5 (repeat of the header, columns C:S only) ---
$ws2.Range("C5").Value2 = "Year"
$ws2.Range("D5").Value2 = " tot in HRUs reaches and reservoirs at end of last year (mm H2O)"
$ws2.Range("E5").Value2 = " Precip (mm H2O)"
$ws2.Range("F5").Value2 = " GW pumping (mm H2O)"
$ws2.Range("G5").Value2 = " High Cascades groundwater contribution mm H2O"
$ws2.Range("H5").Value2 = " from outside the basin (mm H2O)"
$ws2.Range("I5").Value2 = " water added by FlowModel (mm)"
$ws2.Range("J5").Value2 = " to outside the basin (mm H2O)"
$ws2.Range("K5").Value2 = " AET (mm H2O)"
$ws2.Range("L5").Value2 = " SNOW_EVAP (mm H2O)"
$ws2.Range("M5").Value2 = " basin discharge (mm H2O)"
$ws2.Range("N5").Value2 = " tot in HRUs reaches and reservoirs at end of this year (mm H2O)"
$ws2.Range("O5").Value2 = " irrigation (ac-ft)"
$ws2.Range("P5").Value2 = " municipal and rural domestic (ac-ft)"
$ws2.Range("Q5").Value2 = " mass balance discrepancy (mm H2O)"
$ws2.Range("R5").Value2 = " mass balance discrepancy (fraction)"
$ws2.Range("S5").Value2 = " weather year"

$ws2.Range("C5:S5").WrapText = $true
$ws2.Range("D5:N5").NumberFormat = "0.00"
$ws2.Range("O5:P5").NumberFormat = "0"
$ws2.Range("Q5").NumberFormat = "0.00"
$ws2.Range("R5").NumberFormat = "0.000000"

# --- New-sheet view: selected cell E3, not the frozen/scrolled layout sheet1 has ---
$ws2.Range("E3").Select()

# --- sheet1's selection moves now that it is no longer the active sheet ---
$ws1.Range("A1:B1").Select()

$ws2.Activate()
